# Add two new weekly price rows (49 and 50) for Espárragos at
# Terminal Hortofrutícola Agro Chillán, mirroring the existing row layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 49 ---
$ws.Range("A49").Value = 7
$ws.Range("B49").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C49").Value = "Ñuble"
$ws.Range("D49").Value = 45191
$ws.Range("D49").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E49").Value = 16
$ws.Range("F49").Value = 300000000
$ws.Range("G49").Value = "Espárragos"
$ws.Range("H49").Value = "Sin especificar"
$ws.Range("I49").Value = "Primera"
$ws.Range("J49").Value = 200
$ws.Range("K49").Value = 1500
$ws.Range("L49").Value = 1500
$ws.Range("M49").Value = 1500
$ws.Range("N49").Value = "$/kilo"
$ws.Range("O49").Value = "Región del Maule"
$ws.Range("P49").Value = 1500
$ws.Range("Q49").Value = 1
$ws.Range("R49").Value = "Hortaliza"

# --- Row 50 ---
$ws.Range("A50").Value = 7
$ws.Range("B50").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C50").Value = "Ñuble"
$ws.Range("D50").Value = 45191
$ws.Range("D50").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E50").Value = 16
$ws.Range("F50").Value = 300000000
$ws.Range("G50").Value = "Espárragos"
$ws.Range("H50").Value = "Sin especificar"
$ws.Range("I50").Value = "Segunda"
$ws.Range("J50").Value = 200
$ws.Range("K50").Value = 1300
$ws.Range("L50").Value = 1300
$ws.Range("M50").Value = 1300
$ws.Range("N50").Value = "$/kilo"
$ws.Range("O50").Value = "Región del Maule"
$ws.Range("P50").Value = 1300
$ws.Range("Q50").Value = 1
$ws.Range("R50").Value = "Hortaliza"
